# docs/ValueSet-ms-rh-codes.xlsx — "updated docs for pages"
#
# Metadata sheet:
#   B7 (next to "Experimental") goes from blank to the literal text "true"
#   B8 (the "Date" value) is refreshed to a newer timestamp
#
# The "Include ValueSets" / "Include ValueSets 2" sheets are unaffected in
# content (their apparent shared-string index shifts in the source diff are
# just a side effect of the new "true" string being inserted into the
# shared-strings table upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Force this in as literal text (not the Boolean TRUE) by using Excel's
# leading-apostrophe text-entry marker.
$ws.Range("B7").Value = "'true"

$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
